# Applies the "Corrected excel sheets for application fix issues" edit:
#  - Summary sheet: a few corrected figures + new selection
#  - Repayment schedule: re-dated / re-amortised schedule rows, the spare
#    "heading" (O) column's stray data values cleared out, new selection
#  - Transactions: a couple of corrected IDs + new selection
#  - The workbook ends up with "Transactions" as the active/selected tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Select() | Out-Null

$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 700
$wsSummary.Range("E3").Value = 500

$wsSummary.Range("E5").Select() | Out-Null

# ---------------------------------------------------------------------
# Repayment schedule
# ---------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Select() | Out-Null

# The "heading" spacer column (O) should have no data in its body rows,
# same as column E - clear out the stray zeros that were left in O2:O13.
for ($r = 2; $r -le 13; $r++) {
    $wsRepay.Cells.Item($r, 15).Clear() | Out-Null
}
# P2 likewise had a stray empty-but-styled cell that should disappear.
$wsRepay.Range("P2").Clear() | Out-Null

# Row 4 - schedule shifts one period later, interest/due recomputed
$wsRepay.Range("B4").Value = 31
$wsRepay.Range("C4").Value = 42095
$wsRepay.Range("H4").Value = 90.91
$wsRepay.Range("G2").Copy()
$wsRepay.Range("K4").PasteSpecial(-4122) | Out-Null
$wsRepay.Range("K4").Value = 1000
$wsRepay.Range("G2").Copy()
$wsRepay.Range("P4").PasteSpecial(-4122) | Out-Null
$wsRepay.Range("P4").Value = 1000

# Row 5
$wsRepay.Range("B5").Value = 30
$wsRepay.Range("C5").Value = 42125
$wsRepay.Range("H5").Value = 81.82
$wsRepay.Range("H6").Copy()
$wsRepay.Range("K5").PasteSpecial(-4122) | Out-Null
$wsRepay.Range("K5").Value = 990.91
$wsRepay.Range("H6").Copy()
$wsRepay.Range("P5").PasteSpecial(-4122) | Out-Null
$wsRepay.Range("P5").Value = 990.91

# Row 6
$wsRepay.Range("B6").Value = 31
$wsRepay.Range("C6").Value = 42156

# Row 7
$wsRepay.Range("B7").Value = 30
$wsRepay.Range("C7").Value = 42186

# Row 8
$wsRepay.Range("B8").Value = 31
$wsRepay.Range("C8").Value = 42217

# Row 9 (days unchanged)
$wsRepay.Range("C9").Value = 42248

# Row 10
$wsRepay.Range("B10").Value = 30
$wsRepay.Range("C10").Value = 42278

# Row 11
$wsRepay.Range("B11").Value = 31
$wsRepay.Range("C11").Value = 42309

# Row 12
$wsRepay.Range("B12").Value = 30
$wsRepay.Range("C12").Value = 42339

# Row 13
$wsRepay.Range("B13").Value = 31
$wsRepay.Range("C13").Value = 42370

$wsRepay.Range("O11").Select() | Out-Null

# ---------------------------------------------------------------------
# Transactions (ends up the active sheet/tab)
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Select() | Out-Null

$wsTrans.Range("A2").Value = 6372
$wsTrans.Range("A3").Value = 6370

$wsTrans.Range("D3").Select() | Out-Null
